$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 24, shifting existing rows 24:89 down to 25:90.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly record.
$ws.Range("A24").Value = 8
$ws.Range("B24").Value = "Terminal La Palmera de La Serena"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44708
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 100112052
$ws.Range("G24").Value = "Albahaca"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 1120
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 4500
$ws.Range("M24").Value = 4250
$ws.Range("N24").Value = "$/paquete"
$ws.Range("O24").Value = "Región de Arica y Parinacota"
$ws.Range("P24").Value = 4250
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"
